$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header (shared string) renames in row 1 ---
$ws.Range("C1").Value = "GDP"
$ws.Range("E1").Value = "Budget_Previous_Year"
$ws.Range("F1").Value = "LatinAmerica"
$ws.Range("G1").Value = "Africa"
$ws.Range("H1").Value = "Confessional"
$ws.Range("I1").Value = "Universal"
$ws.Range("AF1").Value = "Donor_Aid_Budget"
$ws.Range("AG1").Value = "Total_Funds"
$ws.Range("AH1").Value = "%_Private_Funds"
$ws.Range("AI1").Value = "%_MAE_Funds"
$ws.Range("AM1").Value = "Delegation"

# --- Updated data values (GDP column C, Colony flag column AL) ---
$ws.Range("C2").Value = 5596.139681459835
$ws.Range("C3").Value = 2934.187009790061
$ws.Range("C4").Value = 2870.311589353206
$ws.Range("C5").Value = 1873.394108966653
$ws.Range("C6").Value = 1460.056109840828
$ws.Range("C7").Value = 5191.140356354663
$ws.Range("AL7").Value = 1
$ws.Range("C8").Value = 9502.243585046588
$ws.Range("C9").Value = 17288.8595992193
$ws.Range("C10").Value = 1094.710717322873
$ws.Range("C11").Value = 473.2998774917226
$ws.Range("C12").Value = 5730.354774594881
$ws.Range("C13").Value = 2983.242707849043
$ws.Range("C14").Value = 2898.942214704482
$ws.Range("C15").Value = 665.6274194933962
$ws.Range("AL15").Value = 1
$ws.Range("C16").Value = 1904.346464968814
$ws.Range("C17").Value = 1503.870423231357
$ws.Range("C18").Value = 5555.389721901988
$ws.Range("AL18").Value = 1
$ws.Range("C19").Value = 1955.461557360978
$ws.Range("C20").Value = 492.3430015592067
$ws.Range("C21").Value = 17610.30663334184
$ws.Range("C22").Value = 466.0709276378625
$ws.Range("C23").Value = 5885.254624554112
$ws.Range("C24").Value = 2965.153206179127
$ws.Range("C25").Value = 1939.33862702996
$ws.Range("C26").Value = 1577.487171555845
$ws.Range("C27").Value = 5660.517066940175
$ws.Range("AL27").Value = 1
$ws.Range("C28").Value = 10883.31535948899
$ws.Range("C29").Value = 6255.426161047989
$ws.Range("C30").Value = 3137.260298393558
$ws.Range("C31").Value = 16146.07242861928
$ws.Range("C32").Value = 1775.027517189621
$ws.Range("C33").Value = 515.8271637832048
$ws.Range("C34").Value = 6301.696269820412
$ws.Range("AL34").Value = 1
$ws.Range("C35").Value = 6522.736799041846
$ws.Range("C36").Value = 3210.869677115934
$ws.Range("C37").Value = 14093.81249338665
$ws.Range("C38").Value = 1836.014008604312
$ws.Range("C39").Value = 517.8609592583078
$ws.Range("C40").Value = 6661.86504232374
$ws.Range("AL40").Value = 1
$ws.Range("C41").Value = 7055.001624869326
$ws.Range("C42").Value = 526.5953412037009
